# Auto-generated edit script: updates cryptos list values (prices, 1h volume %, and
# some row reordering for coins whose ranking position swapped) to match the new
# scrape snapshot referenced by the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "87.941.14"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.86%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.110.61"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.18%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.99"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "634.55"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.383"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.96%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.843"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +22.35%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.110.88"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.596"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.24%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.61%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000246"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.37"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "87.901.67"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.684.75"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "31.96"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.148.13"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.37"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000216"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +10.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.29"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "424.45"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.45"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.92"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.50"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +6.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "83.69"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +11.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.43"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.276.38"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.77%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.157"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -5.73%  "
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.15"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.38%  "
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "dogwifhat"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.82"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -8.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.148"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +16.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "501.47"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -5.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.80"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.94%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "PancakeSwap"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.83"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.19%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.27"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "22.53"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.17"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.57%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.34%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.140"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +13.20%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "PolygonEcosystemToken"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.368"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.22%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.84"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "146.61"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.91"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0662"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +12.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "161.69"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -6.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.719"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.19"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.78%  "
